$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: add a new "Utility Panel" row to the parts table + a hyperlink
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$partName = "Utility Panel (Common: 1/8 in. x 4 ft. x 8 ft.; Actual: 0.106 in. x 48 in. x 96 in.)"
$ws1.Range("B4").Value = $partName
$ws1.Range("C4").Value = 96
$ws1.Range("D4").Formula = "=4*12"
$ws1.Range("E4").Value = 0.106
$ws1.Range("F4").Value = 11.44
$ws1.Range("F4").NumberFormat = $ws1.Range("F3").NumberFormat

# Grow Table1 so it covers the new row too
$tbl1 = $ws1.ListObjects.Item(1)
$tbl1.Resize($ws1.Range("A1:F4"))

# Link the new part name to its product page, like the other hyperlinked parts
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://www.homedepot.com/p/1-8-in-x-4-ft-x-8-ft-Utility-Panel/100109169")

# Move the selection like the author left it
$ws1.Activate()
$ws1.Range("E5").Select()

# ---------------------------------------------------------------------------
# Sheet2 (tab "Sheet2", the scratch-pad sheet) - author redid these calcs
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet2")
$ws3.Activate()
$ws3.Cells.ClearContents()

$ws3.Range("A1").Value = "Feet"
$ws3.Range("B1").Value = "Inches"
$ws3.Range("C1").Value = "Inches"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Formula = "=7+23/32"
$ws3.Range("C2").Formula = "=12*A2+B2"
$ws3.Range("D2").Formula = "=C2/2"

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Formula = "=7+11/16"
$ws3.Range("C3").Formula = "=12*A3+B3"
$ws3.Range("D3").Formula = "=C3/2"

$ws3.Range("C4").Formula = "=C2-C3"
$ws3.Range("D4").Formula = "=D2-D3"

$ws3.Range("C5").Formula = "=C4*32"

$ws3.Range("E8").Value = "Top to top of rail"
$ws3.Range("E9").Value = "Rail height"

$ws3.Range("G11").Formula = "=5.5/2"

$ws3.Range("G13").Formula = "=G12+C4/2"
$ws3.Range("G13").NumberFormat = "0.000000"

$ws3.Range("G12").Formula = "=27/32"
$ws3.Range("G12").NumberFormat = "0.00000"

$ws3.Columns.Item(7).AutoFit() | Out-Null

$ws3.Range("G13").Select()

# ---------------------------------------------------------------------------
# Wire_Shelf becomes the active tab (as it was left selected last)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Wire_Shelf")
$ws2.Activate()
